$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update existing values in row 2
$ws.Range("E2").Value = 0.00029628
$ws.Range("F2").Value = 0.01628703
$ws.Range("G2").Value = 0.0005126685623611551

# Update existing values in row 3
$ws.Range("E3").Value = 0.00195624
$ws.Range("F3").Value = 0.01178469
$ws.Range("G3").Value = 0.002833600783645656

# Add new row 4
$ws.Range("A4").Value = "Product_Mode"
$ws.Range("B4").Value = "FEASIBLE_POINT"
$ws.Range("C4").Value = "LOCALLY_SOLVED"
$ws.Range("D4").Value = 0.0
$ws.Range("E4").Value = 0.00394875
$ws.Range("F4").Value = 0.00962001
$ws.Range("G4").Value = 0.005086831985743381
$ws.Range("H4").Value = 5447
$ws.Range("I4").Value = 0.1609039306640625
$ws.Range("J4").Value = 50.81
$ws.Range("K4").Value = 17.8
$ws.Range("L4").Value = 126.75
$ws.Range("M4").Value = 163.24
